$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (shifting existing rows 6-23 down to 7-24)
# for the new "Style of Delivery" relationship (BFO:0000133).
$ws.Rows("6:6").Insert()

# Fill cell values in the same column order/meaning as the pre-existing
# rows: A=ID, B=Relationship, C=Equivalent to relationship, E=Definition.
# Values are assigned in B, C, A, E order so the shared-string table gets
# the new strings appended in the same order as the target workbook.
$ws.Range("B6").Value = "is process attribute of"
$ws.Range("C6").Value = "process profile of [BFO:0000133]"
$ws.Range("A6").Value = "BFO:0000133"
$ws.Range("E6").Value = "inverse of has_process_profile"

# Update the selection/active cell to match the saved view state.
$ws.Range("A6").Select()
